# Applies the "GIANT_HONET" monster-row insertion plus assorted stat tweaks
# to RoStatProcessing.xlsx (StatDef + ClassDef sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# StatDef sheet (monster stat table)
# ---------------------------------------------------------------------
$stat = $wb.Worksheets.Item("StatDef")

# --- isolated single-cell tweaks (unrelated to the row insert below) ---
$stat.Cells.Item(102, 10).Value = 90    # J102  Agi 100 -> 90

$stat.Cells.Item(103, 10).Value = 80    # J103  Agi 100 -> 80
$stat.Cells.Item(103, 12).Value = 115   # L103  Attack 110 -> 115
$stat.Cells.Item(103, 15).Value = 125   # O103  Def 100 -> 125
$stat.Cells.Item(103, 16).Value = 0     # P103  MDef 100 -> 0
$stat.Cells.Item(103, 17).Value = 108   # Q103  Exp 105 -> 108

$stat.Cells.Item(141, 10).Value = 70    # J141  Agi 75 -> 70

# --- insert a new monster row at 302 (Giant Hornet), shifting the rest down ---
$stat.Rows(302).Insert()

# NOTE: new shared-string entries are appended to the workbook's string
# table in first-write order, so the Code / sprite-file / Name values are
# written in that specific sequence to reproduce the canonical ordering
# (GIANT_HONET, giant_honet.spr, Giant Hornet).
$stat.Cells.Item(302, 1).Value = 4299               # Id
$stat.Cells.Item(302, 2).Value = "GIANT_HONET"      # Code
$stat.Cells.Item(302, 33).Value = "giant_honet.spr" # ClientOffset
$stat.Cells.Item(302, 3).Value = "Giant Hornet"     # Name
$stat.Cells.Item(302, 4).Value = 56                 # Level
$stat.Cells.Item(302, 5).Value = 100                # HP
$stat.Cells.Item(302, 6).Value = 100                # Str
$stat.Cells.Item(302, 7).Value = 100                # Int
$stat.Cells.Item(302, 8).Value = 100                # Vit
$stat.Cells.Item(302, 9).Value = 100                # Dex
$stat.Cells.Item(302, 10).Value = 100               # Agi
$stat.Cells.Item(302, 11).Value = 100                # Luk
$stat.Cells.Item(302, 12).Value = 100               # Attack
$stat.Cells.Item(302, 13).Value = 10                # Variance
$stat.Cells.Item(302, 14).Value = 1                 # Range
$stat.Cells.Item(302, 15).Value = 100                # Def
$stat.Cells.Item(302, 16).Value = 100               # MDef
$stat.Cells.Item(302, 17).Value = 100               # Exp
$stat.Cells.Item(302, 18).Value = 100               # JExp
$stat.Cells.Item(302, 19).Value = 10                # ScanDist
$stat.Cells.Item(302, 20).Value = 12                # ChaseDist
$stat.Cells.Item(302, 21).Value = "Small"           # Size
$stat.Cells.Item(302, 22).Value = "Insect"          # Race
$stat.Cells.Item(302, 23).Value = "Wind1"           # Element
$stat.Cells.Item(302, 24).Value = 1292              # RechargeTime
$stat.Cells.Item(302, 25).Value = 340               # HitTime
$stat.Cells.Item(302, 26).Value = 792                # AttackTime
$stat.Cells.Item(302, 27).Value = 155                # MoveSpeed
$stat.Cells.Item(302, 28).Value = "Elite"           # Special
$stat.Cells.Item(302, 29).Value = "Normal"          # Class
$stat.Cells.Item(302, 30).Value = "AiPassive"       # MonsterAiType
$stat.Cells.Item(302, 32).Value = 825               # ClientSprite
$stat.Cells.Item(302, 34).Value = 0                 # ClientShadow
$stat.Cells.Item(302, 35).Value = 0.5               # ClientSize
$stat.Cells.Item(302, 36).Value = 1                 # Flags

# --- grow the Table1 ListObject range by one row to match the new dimension ---
$statTable = $stat.ListObjects.Item(1)
$statTable.Resize($stat.Range("A1:AK414"))

# --- view state: pane / selection, best effort ---
$stat.Application.ActiveWindow.Panes.Item(4).ScrollColumn = 12
$stat.Range("AB302").Select()

# ---------------------------------------------------------------------
# ClassDef sheet
# ---------------------------------------------------------------------
$class = $wb.Worksheets.Item("ClassDef")
$class.Cells.Item(23, 5).Value = 105   # E23  110 -> 105
$class.Cells.Item(24, 5).Value = 105   # E24  115 -> 105
$class.Cells.Item(25, 5).Value = 110   # E25  120 -> 110
$class.Cells.Item(26, 5).Value = 110   # E26  125 -> 110
$class.Cells.Item(27, 5).Value = 115   # E27  130 -> 115
$class.Cells.Item(28, 5).Value = 120   # E28  135 -> 120
$class.Cells.Item(29, 5).Value = 125   # E29  140 -> 125

$class.Range("E25").Select()
